$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-updated data (5 sending clusters x 3 target clusters = 15 rows)
$data = @(
    @("ECs", "Efna4", "Epha7", "ECs", 3, 1, 1.974568, 5.923704, 0.5990695552080697, 0.5990695552080698, 1, 0.3333333333333333, 0.2282403333333334, 0.684721, 0.1514399067192403, 0.1514399067192403, 0.4506760585093333, 4.056084526584, 0.09072303755904687, 0.09072303755904687),
    @("ECs", "Efna4", "Epha7", "FAPs", 3, 1, 1.974568, 5.923704, 0.5990695552080697, 0.5990695552080698, 3, 1, 1.151356333333333, 3.454069, 0.763937263734893, 0.763937263734893, 2.273431372397333, 20.460882351576, 0.4576515567925322, 0.4576515567925323),
    @("ECs", "Efna4", "Epha7", "MuSCs", 3, 1, 1.974568, 5.923704, 0.5990695552080697, 0.5990695552080698, 3, 1, 0.127538, 0.3826140000000001, 0.08462282954586674, 0.08462282954586674, 0.251832453584, 2.266492082256, 0.05069496085649069, 0.0506949608564907),
    @("FAPs", "Efna4", "Epha7", "ECs", 3, 1, 1.020259, 3.060777, 0.309539152527049, 0.309539152527049, 1, 0.3333333333333333, 0.2282403333333334, 0.684721, 0.1514399067192403, 0.1514399067192403, 0.2328642542463334, 2.095778288217, 0.046876580384649, 0.04687658038464899),
    @("FAPs", "Efna4", "Epha7", "FAPs", 3, 1, 1.020259, 3.060777, 0.309539152527049, 0.309539152527049, 3, 1, 1.151356333333333, 3.454069, 0.763937263734893, 0.763937263734893, 1.174681661290333, 10.572134951613, 0.2364684932003315, 0.2364684932003315),
    @("FAPs", "Efna4", "Epha7", "MuSCs", 3, 1, 1.020259, 3.060777, 0.309539152527049, 0.309539152527049, 3, 1, 0.127538, 0.3826140000000001, 0.08462282954586674, 0.08462282954586674, 0.130121792342, 1.171096131078, 0.02619407894206851, 0.02619407894206851),
    @("Inflammatory-Mac", "Efna4", "Epha7", "ECs", 1, 0.3333333333333333, 0.077601, 0.232803, 0.02354357842003994, 0.02354357842003994, 1, 0.3333333333333333, 0.2282403333333334, 0.684721, 0.1514399067192403, 0.1514399067192403, 0.017711678107, 0.159405102963, 0.003565437319767968, 0.003565437319767968),
    @("Inflammatory-Mac", "Efna4", "Epha7", "FAPs", 1, 0.3333333333333333, 0.077601, 0.232803, 0.02354357842003994, 0.02354357842003994, 3, 1, 1.151356333333333, 3.454069, 0.763937263734893, 0.763937263734893, 0.08934640282300001, 0.8041176254070002, 0.01798581687673319, 0.01798581687673319),
    @("Inflammatory-Mac", "Efna4", "Epha7", "MuSCs", 1, 0.3333333333333333, 0.077601, 0.232803, 0.02354357842003994, 0.02354357842003994, 3, 1, 0.127538, 0.3826140000000001, 0.08462282954586674, 0.08462282954586674, 0.009897076338000001, 0.08907368704200001, 0.001992324223538786, 0.001992324223538787),
    @("MuSCs", "Efna4", "Epha7", "ECs", 3, 1, 0.1837383333333333, 0.551215, 0.05574487261247628, 0.05574487261247628, 1, 0.3333333333333333, 0.2282403333333334, 0.684721, 0.1514399067192403, 0.1514399067192403, 0.04193649844611112, 0.377428486015, 0.008441998308509342, 0.008441998308509341),
    @("MuSCs", "Efna4", "Epha7", "FAPs", 3, 1, 0.1837383333333333, 0.551215, 0.05574487261247628, 0.05574487261247628, 3, 1, 1.151356333333333, 3.454069, 0.763937263734893, 0.763937263734893, 0.2115482937594445, 1.903934643835, 0.0425855854508253, 0.0425855854508253),
    @("MuSCs", "Efna4", "Epha7", "MuSCs", 3, 1, 0.1837383333333333, 0.551215, 0.05574487261247628, 0.05574487261247628, 3, 1, 0.127538, 0.3826140000000001, 0.08462282954586674, 0.08462282954586674, 0.02343361955666667, 0.21090257601, 0.004717288853141635, 0.004717288853141635),
    @("Resolving-Mac", "Efna4", "Epha7", "ECs", 1, 0.3333333333333333, 0.03989166666666667, 0.119675, 0.01210284123236505, 0.01210284123236505, 1, 0.3333333333333333, 0.2282403333333334, 0.684721, 0.1514399067192403, 0.1514399067192403, 0.009104887297222223, 0.08194398567500001, 0.001832853147267138, 0.001832853147267138),
    @("Resolving-Mac", "Efna4", "Epha7", "FAPs", 1, 0.3333333333333333, 0.03989166666666667, 0.119675, 0.01210284123236505, 0.01210284123236505, 3, 1, 1.151356333333333, 3.454069, 0.763937263734893, 0.763937263734893, 0.04592952306388889, 0.4133657075750001, 0.009245811414470792, 0.009245811414470794),
    @("Resolving-Mac", "Efna4", "Epha7", "MuSCs", 1, 0.3333333333333333, 0.03989166666666667, 0.119675, 0.01210284123236505, 0.01210284123236505, 3, 1, 0.127538, 0.3826140000000001, 0.08462282954586674, 0.08462282954586674, 0.005087703383333334, 0.04578933045000001, 0.001024176670627115, 0.001024176670627115)
)

# Clear the old data range (previously rows 2-17) then write the new 15 rows (rows 2-16)
$ws.Range("A2:T17").Clear() | Out-Null

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
